# Distance table between the 12 cities listed in 'StedenInformatie'
# (Utrecht, Groningen, Leeuwarden, Assen, Zwolle, Arnhem, Haarlem,
#  Middelburg, Maastricht, Lelystad, Den Haag, Den Bosch), in that order
# for both rows and columns.
$wb = $excel.ActiveWorkbook

# Add a new worksheet named 'distancetable' after the last existing sheet (Costs)
$costsSheet = $wb.Worksheets.Item("Costs")
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $costsSheet)
$newSheet.Name = "distancetable"

$data = New-Object 'object[,]' 12,12
$data[0,0] = 0
$data[0,1] = 57.6
$data[0,2] = 25.760000000000005
$data[0,3] = 47.6
$data[0,4] = 70
$data[0,5] = 46
$data[0,6] = 64
$data[0,7] = 39.760000000000005
$data[0,8] = 63.6
$data[0,9] = 84
$data[0,10] = 26.04
$data[0,11] = 27.400000000000002
$data[1,0] = 57.6
$data[1,1] = 0
$data[1,2] = 81.600000000000009
$data[1,3] = 91.2
$data[1,4] = 12.440000000000001
$data[1,5] = 82.4
$data[1,6] = 32.880000000000003
$data[1,7] = 50
$data[1,8] = 122.80000000000001
$data[1,9] = 132.4
$data[1,10] = 66
$data[1,11] = 30.200000000000003
$data[2,0] = 25.760000000000005
$data[2,1] = 81.600000000000009
$data[2,2] = 0
$data[2,3] = 42.800000000000004
$data[2,4] = 94.4
$data[2,5] = 40.800000000000004
$data[2,6] = 83.2
$data[2,7] = 45.2
$data[2,8] = 50.400000000000006
$data[2,9] = 60
$data[2,10] = 22.400000000000002
$data[2,11] = 51.6
$data[3,0] = 47.6
$data[3,1] = 91.2
$data[3,2] = 42.800000000000004
$data[3,3] = 0
$data[3,4] = 94.4
$data[3,5] = 23.680000000000003
$data[3,6] = 77.2
$data[3,7] = 45.2
$data[3,8] = 92.4
$data[3,9] = 45.6
$data[3,10] = 27.680000000000003
$data[3,11] = 61.6
$data[4,0] = 70
$data[4,1] = 12.440000000000001
$data[4,2] = 94.4
$data[4,3] = 94.4
$data[4,4] = 0
$data[4,5] = 80.800000000000011
$data[4,6] = 25.84
$data[4,7] = 51.2
$data[4,8] = 134.4
$data[4,9] = 142
$data[4,10] = 75.2
$data[4,11] = 42
$data[5,0] = 46
$data[5,1] = 82.4
$data[5,2] = 40.800000000000004
$data[5,3] = 23.680000000000003
$data[5,4] = 80.800000000000011
$data[5,5] = 0
$data[5,6] = 56.400000000000006
$data[5,7] = 32.04
$data[5,8] = 91.2
$data[5,9] = 66.400000000000006
$data[5,10] = 24.040000000000003
$data[5,11] = 53.2
$data[6,0] = 64
$data[6,1] = 32.880000000000003
$data[6,2] = 83.2
$data[6,3] = 77.2
$data[6,4] = 25.84
$data[6,5] = 56.400000000000006
$data[6,6] = 0
$data[6,7] = 39.360000000000007
$data[6,8] = 129.6
$data[6,9] = 130.4
$data[6,10] = 63.2
$data[6,11] = 36.800000000000004
$data[7,0] = 39.760000000000005
$data[7,1] = 50
$data[7,2] = 45.2
$data[7,3] = 45.2
$data[7,4] = 51.2
$data[7,5] = 32.04
$data[7,6] = 39.360000000000007
$data[7,7] = 0
$data[7,8] = 95.2
$data[7,9] = 92.800000000000011
$data[7,10] = 25.64
$data[7,11] = 20.560000000000002
$data[8,0] = 63.6
$data[8,1] = 122.80000000000001
$data[8,2] = 50.400000000000006
$data[8,3] = 92.4
$data[8,4] = 134.4
$data[8,5] = 91.2
$data[8,6] = 129.6
$data[8,7] = 95.2
$data[8,8] = 0
$data[8,9] = 74.8
$data[8,10] = 72
$data[8,11] = 92.4
$data[9,0] = 84
$data[9,1] = 132.4
$data[9,2] = 60
$data[9,3] = 45.6
$data[9,4] = 142
$data[9,5] = 66.400000000000006
$data[9,6] = 130.4
$data[9,7] = 92.800000000000011
$data[9,8] = 74.8
$data[9,9] = 0
$data[9,10] = 68
$data[9,11] = 100.4
$data[10,0] = 26.04
$data[10,1] = 66
$data[10,2] = 22.400000000000002
$data[10,3] = 27.680000000000003
$data[10,4] = 75.2
$data[10,5] = 24.040000000000003
$data[10,6] = 63.2
$data[10,7] = 25.240000000000002
$data[10,8] = 72
$data[10,9] = 68
$data[10,10] = 0
$data[10,11] = 35.92
$data[11,0] = 27.400000000000002
$data[11,1] = 30.200000000000003
$data[11,2] = 51.6
$data[11,3] = 61.6
$data[11,4] = 42
$data[11,5] = 53.2
$data[11,6] = 36.800000000000004
$data[11,7] = 20.560000000000002
$data[11,8] = 92.4
$data[11,9] = 100.4
$data[11,10] = 35.92
$data[11,11] = 0

$newSheet.Range("A1:L12").Value = $data

# Reproduce the saved selection state on the new (now active) sheet
$null = $newSheet.Range("H14").Select()
